$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the "Text" number format (same as A2) to A3:A4 BEFORE writing
# values, so the numeric-looking strings "2"/"3" stay text (shared
# strings) instead of being auto-coerced into numeric cells.
$ws.Range("A3:A4").NumberFormat = "@"

$ws.Range("A3").Value = "2"
$ws.Range("A4").Value = "3"

$ws.Range("B3").Value = "Samsung Galaxy"
$ws.Range("B4").Value = "Iphone"

$ws.Range("C3").Value = "Tecnología"
$ws.Range("C4").Value = "Tecnología"

$ws.Range("C10").Select()
